# Swap columns D<->F and E<->G (codeforiati:category-code/group-name and
# codeforiati:group-code/category-name) across the whole used range of the
# active worksheet. Column A..C are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 4).End(-4162).Row  # xlUp = -4162
if ($lastRow -lt 1) { $lastRow = 1 }

for ($r = 1; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)
    $gCell = $ws.Cells.Item($r, 7)

    $dVal = $dCell.Text
    $eVal = $eCell.Text
    $fVal = $fCell.Text
    $gVal = $gCell.Text

    # Force text (string) storage - these columns hold text-typed codes
    # (e.g. "110") that Excel would otherwise coerce to numbers.
    $dCell.NumberFormat = "@"
    $eCell.NumberFormat = "@"
    $fCell.NumberFormat = "@"
    $gCell.NumberFormat = "@"

    $dCell.Value = "'" + $fVal
    $eCell.Value = "'" + $gVal
    $fCell.Value = "'" + $eVal
    $gCell.Value = "'" + $dVal
}
